$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: scenario code DKW -> SHR
$ws.Range("G3").Value = "SHR"

# D4:D11 used to hold the literal number 2012; they now hold the text "2012-0"
# (same format/style as the "2012-0" label already used lower on the sheet).
$ws.Range("D4:D11").NumberFormat = "0.00"
$ws.Range("D4:D11").Value = "2012-0"

# D15:D22 previously repeated the "2012-0" label; it has been removed entirely.
$ws.Range("D15:D22").Clear()

# Move the active selection from D4:D11 to J3.
$ws.Range("J3").Select()
